# Regenerate save_data: update column G ("K") values for rows 2-25
# to reflect strike count (K) computed from the regenerated std/mean
# and s_vals, instead of the old "Strike#" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 1
    6  = 0
    7  = 2
    8  = 0
    9  = 2
    10 = 1
    11 = 2
    12 = 2
    13 = 0
    14 = 3
    15 = 0
    16 = 4
    17 = 1
    18 = 1
    19 = 2
    20 = 2
    21 = 0
    22 = 2
    23 = 3
    24 = 1
    25 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
